# Update computed price/profit columns (H-N) per the scheduled market-data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 10500
$ws.Range("I32").Value = 6500
$ws.Range("J32").Value = 14500
$ws.Range("K32").Value = 6500
$ws.Range("L32").Value = 14500
$ws.Range("M32").Value = -6174
$ws.Range("N32").Value = -15152
$ws.Range("H55").Value = 150
$ws.Range("I55").Value = 150
$ws.Range("K55").Value = 150
$ws.Range("M55").Value = 64
$ws.Range("H137").Value = 2484.3823
$ws.Range("I137").Value = 1512.3158
$ws.Range("J137").Value = 3715.6667
$ws.Range("K137").Value = 4536.9474
$ws.Range("L137").Value = 11147.0001
$ws.Range("M137").Value = -1986.9474
$ws.Range("N137").Value = -16247.0001
$ws.Range("H138").Value = 4353.676
$ws.Range("J138").Value = 4735.1
$ws.Range("L138").Value = 14205.3
$ws.Range("N138").Value = -24485.3

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10440.206
$ws.Range("I32").Value = 10131.322
$ws.Range("K32").Value = 10131.322
$ws.Range("M32").Value = -9844.322
$ws.Range("H44").Value = 40000
$ws.Range("J44").Value = 40000
$ws.Range("L44").Value = 40000
$ws.Range("N44").Value = -40976
$ws.Range("H45").Value = 3212.4285
$ws.Range("J45").Value = 3578.25
$ws.Range("L45").Value = 3578.25
$ws.Range("N45").Value = -4332.25
$ws.Range("H132").Value = 1795.9032
$ws.Range("I132").Value = 1157.32
$ws.Range("J132").Value = 4456.6665
$ws.Range("K132").Value = 3471.96
$ws.Range("L132").Value = 13369.9995
$ws.Range("M132").Value = -941.96
$ws.Range("N132").Value = -18429.9995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 1239.75
$ws.Range("J64").Value = 1319.6666
$ws.Range("L64").Value = 1319.6666
$ws.Range("N64").Value = -1769.6666
$ws.Range("H67").Value = 1239.75
$ws.Range("J67").Value = 1319.6666
$ws.Range("L67").Value = 1319.6666
$ws.Range("N67").Value = -2879.6666
$ws.Range("H107").Value = 1069.5333
$ws.Range("I107").Value = 926.53845
$ws.Range("K107").Value = 926.53845
$ws.Range("M107").Value = 993.46155

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2607.4546
$ws.Range("I31").Value = 2630
$ws.Range("J31").Value = 2506
$ws.Range("K31").Value = 2630
$ws.Range("L31").Value = 2506
$ws.Range("M31").Value = -2335
$ws.Range("N31").Value = -3096
$ws.Range("H34").Value = 2607.4546
$ws.Range("I34").Value = 2630
$ws.Range("J34").Value = 2506
$ws.Range("K34").Value = 2630
$ws.Range("L34").Value = 2506
$ws.Range("M34").Value = -2428
$ws.Range("N34").Value = -2910
$ws.Range("H86").Value = 22027.4
$ws.Range("I86").Value = 10021.667
$ws.Range("J86").Value = 40036
$ws.Range("K86").Value = 10021.667
$ws.Range("L86").Value = 40036
$ws.Range("M86").Value = -8898.666999999999
$ws.Range("N86").Value = -42282
$ws.Range("H89").Value = 22027.4
$ws.Range("I89").Value = 10021.667
$ws.Range("J89").Value = 40036
$ws.Range("K89").Value = 50108.335
$ws.Range("L89").Value = 200180
$ws.Range("M89").Value = -44492.335
$ws.Range("N89").Value = -211412

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 302.57144
$ws.Range("J23").Value = 302.57144
$ws.Range("L23").Value = 907.71432
$ws.Range("N23").Value = -1377.71432
$ws.Range("H46").Value = 3266.6667
$ws.Range("I46").Value = 2900
$ws.Range("J46").Value = 3450
$ws.Range("K46").Value = 8700
$ws.Range("L46").Value = 10350
$ws.Range("M46").Value = -8609
$ws.Range("N46").Value = -10532
$ws.Range("H128").Value = 243333
$ws.Range("I128").Value = 243333
$ws.Range("K128").Value = 729999
$ws.Range("M128").Value = -725019
$ws.Range("H138").Value = 6799.5
$ws.Range("I138").Value = 6799.5
$ws.Range("K138").Value = 20398.5
$ws.Range("M138").Value = -15258.5
$ws.Range("H140").Value = 1348
$ws.Range("I140").Value = 1348
$ws.Range("K140").Value = 4044
$ws.Range("M140").Value = 1136
$ws.Range("H141").Value = 5030.3335
$ws.Range("I141").Value = 5030.3335
$ws.Range("K141").Value = 15091.0005
$ws.Range("M141").Value = -9911.000499999998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("M29").ClearContents()
$ws.Range("N29").ClearContents()
$ws.Range("H80").Value = 5987.778
$ws.Range("I80").Value = 4699.6
$ws.Range("J80").Value = 7598
$ws.Range("K80").Value = 4699.6
$ws.Range("L80").Value = 7598
$ws.Range("M80").Value = -3701.6
$ws.Range("N80").Value = -9594
$ws.Range("H83").Value = 5987.778
$ws.Range("I83").Value = 4699.6
$ws.Range("J83").Value = 7598
$ws.Range("K83").Value = 23498
$ws.Range("L83").Value = 37990
$ws.Range("M83").Value = -18506
$ws.Range("N83").Value = -47974
$ws.Range("H97").Value = 584
$ws.Range("I97").Value = 484.58823
$ws.Range("K97").Value = 484.58823
$ws.Range("M97").Value = 11.41176999999999
$ws.Range("H102").Value = 2949.6667
$ws.Range("I102").Value = 2924.75
$ws.Range("K102").Value = 2924.75
$ws.Range("M102").Value = -1302.75
$ws.Range("H126").Value = 5045.875
$ws.Range("I126").Value = 3978.5
$ws.Range("K126").Value = 11935.5
$ws.Range("M126").Value = -9465.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 966
$ws.Range("J46").Value = 2000
$ws.Range("L46").Value = 2000
$ws.Range("N46").Value = -2376
$ws.Range("H55").Value = 375.77777
$ws.Range("I55").Value = 370.66666
$ws.Range("K55").Value = 370.66666
$ws.Range("M55").Value = -197.66666
$ws.Range("H93").Value = 2584.7144
$ws.Range("I93").Value = 2584.7144
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 2584.7144
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = -1336.7144
$ws.Range("N93").ClearContents()
$ws.Range("H132").Value = 5720.091
$ws.Range("I132").Value = 5658.1113
$ws.Range("K132").Value = 16974.3339
$ws.Range("M132").Value = -14444.3339
